$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17609965801239
$ws.Range("B1").Value = 2.411120653152466
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.341566562652588
$ws.Range("E1").Value = 1.205419898033142
